# The 21 "NNN. Problem Name" entries that used to sit right after
# "String/Array" (rows 8-28 of the shared-string table / B71:B91 on the
# sheet) are renumbered/reshuffled upstream; here they keep their row
# positions (B71:B91) but lose their leading "NNN. " numeric prefix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B71").Value = "Populating Next Right Pointers in Each Node II"
$ws.Range("B72").Value = "Text Justification"
$ws.Range("B73").Value = "Maximal Square"
$ws.Range("B74").Value = "H-Index"
$ws.Range("B75").Value = "Maximal Rectangle"
$ws.Range("B76").Value = "Course Schedule II"
$ws.Range("B77").Value = "Reverse Nodes in k-Group"
$ws.Range("B78").Value = "Palindrom Linked List"
$ws.Range("B79").Value = "Task Scheduler"
$ws.Range("B80").Value = "Subtree of Another Tree"
$ws.Range("B81").Value = "Brick Wall"
$ws.Range("B82").Value = "Diameter of Binary Tree"
$ws.Range("B83").Value = "Contiguous Array"
$ws.Range("B84").Value = "Continuous Subarray Sum"
$ws.Range("B85").Value = "Total Hamming Distance"
$ws.Range("B86").Value = "Hamming Distance"
$ws.Range("B87").Value = "Sum of Left Leaves"
$ws.Range("B88").Value = "H-index II"
$ws.Range("B89").Value = "Graph Valid Tree"
$ws.Range("B90").Value = "Remove Duplicates from Sorted Array II"
$ws.Range("B91").Value = "Remove Duplicates from Sorted Array"

# Update the view state to match: scrolled further down, with B91 selected.
$ws.Activate()
$ws.Range("B91").Select()
$excel.ActiveWindow.ScrollRow = 70
$excel.ActiveWindow.ScrollColumn = 1

